$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newName = "1061-MS-EPP-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-LateRepayment"

# Update the product name value (B1) on both sheets, and bold it (style 5).
$wsInput.Range("B1").Value = $newName
$wsInput.Range("B1").Font.Bold = $true

$wsOutput.Range("B1").Value = $newName
$wsOutput.Range("B1").Font.Bold = $true

# Reset selection / scroll position on the input sheet.
$wsInput.Range("B1").Select()

# Make the output sheet the active sheet/tab, with B1 selected.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
